$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8569129109382629
$ws.Range("B1").Value = 1.526379823684692
$ws.Range("C1").Value = 6.273950099945068
$ws.Range("D1").Value = 2.952908039093018
$ws.Range("E1").Value = 1.588418126106262
